$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark that currently sits between
#    "...llevarse a cabo" and the trailing "." in paragraph 4, then
#    rewrite that whole paragraph as a single clean run (merging the
#    "." back into the sentence).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$p4 = $d.Paragraphs.Item(4)
$p4Start = $p4.Range.Start
$p4End = $p4.Range.End
$p4Range = $d.Range($p4Start, $p4End)
$p4Range.Text = "La estimación de los costos de un plan fue una parte bastante interesante pues  es algo completamente necesario, y sin embargo no es nada fácil de realizar. La estimación requiere experiencia para llevarse a cabo."

# ------------------------------------------------------------------
# 2. Re-create the "_GoBack" bookmark at the end of paragraph 1
#    ("Lecciones Aprendidas"). A collapsed, zero-length Range placed
#    exactly at a run/paragraph-mark boundary does not reliably anchor
#    a new bookmark, so insert a unique placeholder run right after the
#    text, use Find to get a real (non-empty) Range around it, anchor
#    the bookmark there, then delete the placeholder text again -
#    leaving bookmarkStart/bookmarkEnd adjacent, right after the run.
# ------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("Lecciones Aprendidas", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRange.Collapse(0)
$findRange.InsertAfter("@@BOOKMARKPLACEHOLDER@@")

$bmRange = $d.Content
$bmRange.Find.Execute("@@BOOKMARKPLACEHOLDER@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange)
$bmRange.Text = ""

# ------------------------------------------------------------------
# 3. Justify (both) the first four paragraphs.
# ------------------------------------------------------------------
for ($i = 1; $i -le 4; $i++) {
    $d.Paragraphs.Item($i).Alignment = 3
}

# ------------------------------------------------------------------
# 4. Append three new justified paragraphs after the (previously
#    empty) trailing paragraph, reusing it for the first one.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastRange.InsertAfter("Durante esta experiencia educativa también aprendí la importancia de los acuerdos a los que se llega con el cliente, que en este caso fue el profesor, pues cuando se llega a un acuerdo con el cliente en el mundo real, lo que se dijo se debe de cumplir. ")
$lastPara.Alignment = 3

$lastRange.Collapse(0)
$lastRange.InsertParagraphAfter()
$lastRange.Collapse(0)
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.InsertAfter("Durante nuestro trabajo en equipo decidimos usar la metodología SCRUM para el desarrollo del proyecto, pues vimos que cuando se aplica esta metodología la documentación se reduce en gran medida. El inconveniente de usar esta metodología es que se requiere gran experiencia de parte de todo el equipo, también aprendí que la documentación es muy necesaria durante el desarrollo de software pues sin ella partes del desarrollo se dificultan. ")
$newPara.Alignment = 3

$newRange = $newPara.Range
$newRange.Collapse(0)
$newRange.InsertParagraphAfter()
$newRange.Collapse(0)
$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara2.Range.InsertAfter("El rol que desempeñe como líder de proyecto durante el curso me enseño que dirigir un equipo de trabajo no es nada fácil, pues se debe tener los conocimientos necesarios de los demás roles del equipo pues si presenta algún problema que un compañero no pueda resolver, el líder debe de saber qué hacer. Por lo tanto el líder de proyecto es un rol para el cual debe de estar muy preparado, tanto en conocimientos técnicos como en tener un trato con las personas. ")
$newPara2.Alignment = 3

Write-Output "done"
